# Auto-generated Excel COM-interop script to apply numeric corrections
# to the Ultros_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR).
# Each row below mirrors one <c> cell change from the source XML diff:
# currentAveragePrice/.../LeveProfit columns (H-N) were recalculated
# from refreshed market data by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3300.6667
$ws.Range("I18").Value = 1267
$ws.Range("J18").Value = 5334.3335
$ws.Range("K18").Value = 1267
$ws.Range("L18").Value = 5334.3335
$ws.Range("M18").Value = -983
$ws.Range("N18").Value = -5902.3335
$ws.Range("H138").Value = 2879.22
$ws.Range("J138").Value = 3018.2273
$ws.Range("L138").Value = 9054.6819
$ws.Range("N138").Value = -19334.6819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19243994
$ws.Range("I32").Value = 20418994
$ws.Range("K32").Value = 20418994
$ws.Range("M32").Value = -20418707
$ws.Range("H74").Value = 1648.75
$ws.Range("I74").Value = 1239.1
$ws.Range("K74").Value = 1239.1
$ws.Range("M74").Value = -365.0999999999999
$ws.Range("H77").Value = 1648.75
$ws.Range("I77").Value = 1239.1
$ws.Range("K77").Value = 6195.5
$ws.Range("M77").Value = -1827.5
$ws.Range("H97").Value = 1007.9167
$ws.Range("I97").Value = 1047.3914
$ws.Range("K97").Value = 1047.3914
$ws.Range("M97").Value = -551.3914
$ws.Range("H110").Value = 1786.8948
$ws.Range("I110").Value = 1830.6111
$ws.Range("K110").Value = 1830.6111
$ws.Range("M110").Value = 214.3888999999999
$ws.Range("H132").Value = 5419.6875
$ws.Range("I132").Value = 4773.6
$ws.Range("K132").Value = 14320.8
$ws.Range("M132").Value = -11790.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 37916.668
$ws.Range("J63").Value = 37916.668
$ws.Range("L63").Value = 37916.668
$ws.Range("N63").Value = -39288.668
$ws.Range("H66").Value = 37916.668
$ws.Range("J66").Value = 37916.668
$ws.Range("L66").Value = 113750.004
$ws.Range("N66").Value = -120614.004
$ws.Range("H107").Value = 1395.2693
$ws.Range("I107").Value = 1223.7
$ws.Range("J107").Value = 1967.1666
$ws.Range("K107").Value = 1223.7
$ws.Range("L107").Value = 1967.1666
$ws.Range("M107").Value = 696.3
$ws.Range("N107").Value = -5807.1666
$ws.Range("H134").Value = 3062.5293
$ws.Range("I134").Value = 1854.3077
$ws.Range("J134").Value = 6989.25
$ws.Range("K134").Value = 5562.9231
$ws.Range("L134").Value = 20967.75
$ws.Range("M134").Value = -3027.9231
$ws.Range("N134").Value = -26037.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49833.418
$ws.Range("J20").Value = 49833.418
$ws.Range("L20").Value = 49833.418
$ws.Range("N20").Value = -50305.418
$ws.Range("H30").Value = 49833.418
$ws.Range("J30").Value = 49833.418
$ws.Range("L30").Value = 49833.418
$ws.Range("N30").Value = -50015.418
$ws.Range("H31").Value = 2766.8164
$ws.Range("I31").Value = 2074.634
$ws.Range("J31").Value = 6314.25
$ws.Range("K31").Value = 2074.634
$ws.Range("L31").Value = 6314.25
$ws.Range("M31").Value = -1779.634
$ws.Range("N31").Value = -6904.25
$ws.Range("H34").Value = 2766.8164
$ws.Range("I34").Value = 2074.634
$ws.Range("J34").Value = 6314.25
$ws.Range("K34").Value = 2074.634
$ws.Range("L34").Value = 6314.25
$ws.Range("M34").Value = -1872.634
$ws.Range("N34").Value = -6718.25
$ws.Range("H99").Value = 8981233
$ws.Range("I99").Value = 1631769.4
$ws.Range("J99").Value = 20005430
$ws.Range("K99").Value = 1631769.4
$ws.Range("L99").Value = 20005430
$ws.Range("M99").Value = -1630271.4
$ws.Range("N99").Value = -20008426
$ws.Range("H126").Value = 8981233
$ws.Range("I126").Value = 1631769.4
$ws.Range("J126").Value = 20005430
$ws.Range("K126").Value = 4895308.199999999
$ws.Range("L126").Value = 60016290
$ws.Range("M126").Value = -4892838.199999999
$ws.Range("N126").Value = -60021230
$ws.Range("H128").Value = 49833.418
$ws.Range("J128").Value = 49833.418
$ws.Range("L128").Value = 49833.418
$ws.Range("N128").Value = -59793.418
$ws.Range("H134").Value = 5041.0527
$ws.Range("I134").Value = 2978.3
$ws.Range("J134").Value = 7333
$ws.Range("K134").Value = 8934.900000000001
$ws.Range("L134").Value = 21999
$ws.Range("M134").Value = -6399.900000000001
$ws.Range("N134").Value = -27069

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7547.7144
$ws.Range("I5").Value = 510.2
$ws.Range("J5").Value = 9746.9375
$ws.Range("K5").Value = 1530.6
$ws.Range("L5").Value = 29240.8125
$ws.Range("M5").Value = -1418.6
$ws.Range("N5").Value = -29464.8125
$ws.Range("H6").Value = 1054.75
$ws.Range("J6").Value = 4
$ws.Range("L6").Value = 12
$ws.Range("N6").Value = -238
$ws.Range("H135").Value = 7547.7144
$ws.Range("I135").Value = 510.2
$ws.Range("J135").Value = 9746.9375
$ws.Range("K135").Value = 4591.8
$ws.Range("L135").Value = 87722.4375
$ws.Range("M135").Value = -2056.8
$ws.Range("N135").Value = -92792.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 27873804
$ws.Range("J80").Value = 55559224
$ws.Range("L80").Value = 55559224
$ws.Range("N80").Value = -55561220
$ws.Range("H83").Value = 27873804
$ws.Range("J83").Value = 55559224
$ws.Range("L83").Value = 277796120
$ws.Range("N83").Value = -277806104
$ws.Range("H97").Value = 11364046
$ws.Range("J97").Value = 50000510
$ws.Range("L97").Value = 50000510
$ws.Range("N97").Value = -50001502
$ws.Range("H132").Value = 4709.4287
$ws.Range("I132").Value = 4680
$ws.Range("J132").Value = 4748.6665
$ws.Range("K132").Value = 14040
$ws.Range("L132").Value = 14245.9995
$ws.Range("M132").Value = -11510
$ws.Range("N132").Value = -19305.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 11758.333
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 11758.333
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 11758.333
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -13630.333
$ws.Range("H77").Value = 11758.333
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 11758.333
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 35274.999
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -44634.999
$ws.Range("H107").Value = 752
$ws.Range("I107").Value = 529.375
$ws.Range("J107").Value = 1197.25
$ws.Range("K107").Value = 1588.125
$ws.Range("L107").Value = 3591.75
$ws.Range("M107").Value = 331.875
$ws.Range("N107").Value = -7431.75
$ws.Range("H113").Value = 1064
$ws.Range("J113").Value = 2020
$ws.Range("L113").Value = 6060
$ws.Range("N113").Value = -10400
$ws.Range("H132").Value = 3042.318
$ws.Range("I132").Value = 2696.5
$ws.Range("K132").Value = 8089.5
$ws.Range("M132").Value = -5559.5
$ws.Range("H135").Value = 77400
$ws.Range("J135").Value = 77400
$ws.Range("L135").Value = 77400
$ws.Range("N135").Value = -87540
